# Scheduled-runner data refresh: update cached market-price / profit
# columns (H:N) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR Leve-profit sheets.
# Values mirror a fresh pull from the price API; some rows gain/lose the
# NQ-profit (M) or HQ-profit (N) cell depending on which price applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 153.5
$ws.Range("I5").Value = 78
$ws.Range("J5").Value = 178.66667
$ws.Range("K5").Value = 78
$ws.Range("L5").Value = 178.66667
$ws.Range("M5").Value = 37
$ws.Range("N5").Value = -408.66667

$ws.Range("N18").ClearContents()
$ws.Range("H18").Value = 625.125
$ws.Range("I18").Value = 625.125
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 625.125
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -341.125

$ws.Range("H40").Value = 1284
$ws.Range("I40").Value = 1260.8695
$ws.Range("J40").Value = 1550
$ws.Range("K40").Value = 1260.8695
$ws.Range("L40").Value = 1550
$ws.Range("M40").Value = -1085.8695
$ws.Range("N40").Value = -1900

$ws.Range("N70").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0

$ws.Range("N73").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0

$ws.Range("H74").Value = 3149.6191
$ws.Range("I74").Value = 3069.2942
$ws.Range("J74").Value = 3491
$ws.Range("K74").Value = 3069.2942
$ws.Range("L74").Value = 3491
$ws.Range("M74").Value = -2133.2942
$ws.Range("N74").Value = -5363

$ws.Range("H77").Value = 3149.6191
$ws.Range("I77").Value = 3069.2942
$ws.Range("J77").Value = 3491
$ws.Range("K77").Value = 15346.471
$ws.Range("L77").Value = 17455
$ws.Range("M77").Value = -10666.471
$ws.Range("N77").Value = -26815

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 392.16666
$ws.Range("I5").Value = 417.75
$ws.Range("J5").Value = 341
$ws.Range("K5").Value = 417.75
$ws.Range("L5").Value = 341
$ws.Range("M5").Value = -305.75
$ws.Range("N5").Value = -565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 392.16666
$ws.Range("I4").Value = 417.75
$ws.Range("J4").Value = 341
$ws.Range("K4").Value = 417.75
$ws.Range("L4").Value = 341
$ws.Range("M4").Value = -302.75
$ws.Range("N4").Value = -571

$ws.Range("H86").Value = 3828.476
$ws.Range("I86").Value = 3186.2856
$ws.Range("J86").Value = 5112.857
$ws.Range("K86").Value = 3186.2856
$ws.Range("L86").Value = 5112.857
$ws.Range("M86").Value = -2063.2856
$ws.Range("N86").Value = -7358.857

$ws.Range("H89").Value = 3828.476
$ws.Range("I89").Value = 3186.2856
$ws.Range("J89").Value = 5112.857
$ws.Range("K89").Value = 15931.428
$ws.Range("L89").Value = 25564.285
$ws.Range("M89").Value = -10315.428
$ws.Range("N89").Value = -36796.285

$ws.Range("H107").Value = 1772.1364
$ws.Range("I107").Value = 1353.3529
$ws.Range("J107").Value = 3196
$ws.Range("K107").Value = 1353.3529
$ws.Range("L107").Value = 3196
$ws.Range("M107").Value = 566.6470999999999
$ws.Range("N107").Value = -7036

$ws.Range("N119").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126
$ws.Range("I7").Value = 93.666664
$ws.Range("J7").Value = 174.5
$ws.Range("K7").Value = 93.666664
$ws.Range("L7").Value = 174.5
$ws.Range("M7").Value = 19.333336
$ws.Range("N7").Value = -400.5

$ws.Range("H22").Value = 728.3214
$ws.Range("I22").Value = 1073.9375
$ws.Range("J22").Value = 267.5
$ws.Range("K22").Value = 1073.9375
$ws.Range("L22").Value = 267.5
$ws.Range("M22").Value = -723.9375
$ws.Range("N22").Value = -967.5

$ws.Range("H31").Value = 3829.1562
$ws.Range("I31").Value = 3621.9443
$ws.Range("J31").Value = 4095.5715
$ws.Range("K31").Value = 3621.9443
$ws.Range("L31").Value = 4095.5715
$ws.Range("M31").Value = -3326.9443
$ws.Range("N31").Value = -4685.5715

$ws.Range("H34").Value = 3829.1562
$ws.Range("I34").Value = 3621.9443
$ws.Range("J34").Value = 4095.5715
$ws.Range("K34").Value = 3621.9443
$ws.Range("L34").Value = 4095.5715
$ws.Range("M34").Value = -3419.9443
$ws.Range("N34").Value = -4499.5715

$ws.Range("H62").Value = 2200
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1576

$ws.Range("H65").Value = 2200
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -7880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M53").ClearContents()
$ws.Range("H53").Value = 19830.75
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 19830.75
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 19830.75
$ws.Range("N53").Value = -21092.75

$ws.Range("H113").Value = 1574
$ws.Range("I113").Value = 1588.8
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1588.8
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 581.2
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 164290.5
$ws.Range("I46").Value = 1286.125
$ws.Range("J46").Value = 245792.69
$ws.Range("K46").Value = 1286.125
$ws.Range("L46").Value = 245792.69
$ws.Range("M46").Value = -1098.125
$ws.Range("N46").Value = -246168.69

$ws.Range("H57").Value = 16469.908
$ws.Range("I57").Value = 1010.25
$ws.Range("J57").Value = 25304
$ws.Range("K57").Value = 1010.25
$ws.Range("L57").Value = 25304
$ws.Range("M57").Value = -444.25
$ws.Range("N57").Value = -26436

$ws.Range("H59").Value = 28000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 28000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 28000
$ws.Range("N59").Value = -29308

$ws.Range("H68").Value = 2895.3157
$ws.Range("I68").Value = 1111
$ws.Range("J68").Value = 2994.4443
$ws.Range("K68").Value = 1111
$ws.Range("L68").Value = 2994.4443
$ws.Range("M68").Value = -362
$ws.Range("N68").Value = -4492.4443

$ws.Range("H71").Value = 2895.3157
$ws.Range("I71").Value = 1111
$ws.Range("J71").Value = 2994.4443
$ws.Range("K71").Value = 5555
$ws.Range("L71").Value = 14972.2215
$ws.Range("M71").Value = -1811
$ws.Range("N71").Value = -22460.2215

$ws.Range("H111").Value = 34266.668
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 34266.668
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 34266.668
$ws.Range("N111").Value = -42446.668

$ws.Range("H132").Value = 6357.1816
$ws.Range("I132").Value = 2779.2104
$ws.Range("J132").Value = 11213
$ws.Range("K132").Value = 8337.6312
$ws.Range("L132").Value = 33639
$ws.Range("M132").Value = -5807.6312
$ws.Range("N132").Value = -38699

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N108").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0

$ws.Range("H114").Value = 62950
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 62950
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 62950
$ws.Range("N114").Value = -71628
